$wb = $excel.ActiveWorkbook

# --- Update Hoja1!A1 text with new conversion rates ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 5.32 = 21260.75 pesos`n✅ 21260.75 pesos = 5.3 = 957.06 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$ws1.Range("A1").Value = $newText

# --- Update tasas sheet N10/O10/N12/O12 values ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 187.999
$ws2.Range("O10").Value = 3997
$ws2.Range("N12").Value = 4012
$ws2.Range("O12").Value = 180.601
